$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff": the localization-status report is
# regenerated. The two tracked files (ad8f1733-... and bf04067f-...) swap
# places in every sheet's first data row vs second data row, bf04067f is
# now the one "Handed back: in sync with en-US", and ad8f1733 is now
# "Ready for handoff" (its handback file is stale vs the latest commit),
# with fresh handoff/handback timestamps and a new error detail message.
# ---------------------------------------------------------------------------

$adName        = "ad8f1733-723b-48bf-9ac9-c27b028504c0.md"
$adPath        = "e2e\ad8f1733-723b-48bf-9ac9-c27b028504c0.md"
$bfName        = "bf04067f-7dea-4840-9a40-b2b7d22fc676.md"
$bfPath        = "e2e\bf04067f-7dea-4840-9a40-b2b7d22fc676.md"

$adUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751fb137d464f3f7d6369474422d04da3a998542/e2e/ad8f1733-723b-48bf-9ac9-c27b028504c0.md"
$bfUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751fb137d464f3f7d6369474422d04da3a998542/e2e/bf04067f-7dea-4840-9a40-b2b7d22fc676.md"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady      = "Ready for handoff"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751fb137d464f3f7d6369474422d04da3a998542/e2e/ad8f1733-723b-48bf-9ac9-c27b028504c0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/982a3946b999b5f64e914594f046542ee42e446c/e2e/ad8f1733-723b-48bf-9ac9-c27b028504c0.md."

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A2").Value = $bfName
$wsOverview.Range("B2").Value = $bfPath
$wsOverview.Range("A3").Value = $adName
$wsOverview.Range("B3").Value = $adPath

$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-08-23 10:49:48"

# Rebuild the hyperlinks so the displayed text follows the swapped rows
# (Address stays pinned to the same github blob per row/file).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $bfUrl, "", "", $bfPath)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $adUrl, "", "", $adPath)

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$wsZh = $wb.Worksheets.Item(2)

$zhBfXlf = "bf04067f-7dea-4840-9a40-b2b7d22fc676.31a53b077f4b72b4183cdb2cba56feb6a2436cbd.zh-cn.xlf"
$zhAdXlf = "ad8f1733-723b-48bf-9ac9-c27b028504c0.0311717db0bcf1e191b5229263166a39ade26397.zh-cn.xlf"

$wsZh.Range("A2").Value = $bfName
$wsZh.Range("G2").Value = $zhBfXlf
$wsZh.Range("I2").Value = $bfName
$wsZh.Range("J2").Value = $zhBfXlf

$wsZh.Range("A3").Value = $adName
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("G3").Value = $zhAdXlf
$wsZh.Range("H3").Value = "2016-08-23 10:49:44"
$wsZh.Range("I3").Value = $adName
$wsZh.Range("J3").Value = $zhAdXlf
$wsZh.Range("P3").Value = $errorDetail

$zhAdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/29a1830b4f09f59bf07533c09e2d89b76b97878b/e2e/ad8f1733-723b-48bf-9ac9-c27b028504c0.md"
$zhBfUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/29a1830b4f09f59bf07533c09e2d89b76b97878b/e2e/bf04067f-7dea-4840-9a40-b2b7d22fc676.md"

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $bfUrl, "", "", $bfName)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhBfUrl, "", "", $bfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $adUrl, "", "", $adName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhAdUrl, "", "", $adName)

# Widen the Error Detail column so the new message is readable.
$zhColWidth = 40 - (5/6)
$wsZh.Columns.Item(16).ColumnWidth = $zhColWidth

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$wsDe = $wb.Worksheets.Item(3)

$deBfXlf = "bf04067f-7dea-4840-9a40-b2b7d22fc676.31a53b077f4b72b4183cdb2cba56feb6a2436cbd.de-de.xlf"
$deAdXlf = "ad8f1733-723b-48bf-9ac9-c27b028504c0.0311717db0bcf1e191b5229263166a39ade26397.de-de.xlf"

$wsDe.Range("A2").Value = $bfName
$wsDe.Range("G2").Value = $deBfXlf
$wsDe.Range("I2").Value = $bfName
$wsDe.Range("J2").Value = $deBfXlf

$wsDe.Range("A3").Value = $adName
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("G3").Value = $deAdXlf
$wsDe.Range("H3").Value = "2016-08-23 10:49:48"
$wsDe.Range("I3").Value = $adName
$wsDe.Range("J3").Value = $deAdXlf
$wsDe.Range("P3").Value = $errorDetail

$deAdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7d4dd02f7889a62e1413c7d40af6d93596f69868/e2e/ad8f1733-723b-48bf-9ac9-c27b028504c0.md"
$deBfUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7d4dd02f7889a62e1413c7d40af6d93596f69868/e2e/bf04067f-7dea-4840-9a40-b2b7d22fc676.md"

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $bfUrl, "", "", $bfName)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deBfUrl, "", "", $bfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $adUrl, "", "", $adName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deAdUrl, "", "", $adName)

# Widen the Error Detail column so the new message is readable.
$deColWidth = 40 - (5/6)
$wsDe.Columns.Item(16).ColumnWidth = $deColWidth
